$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 8.280371333333333
$ws.Range("H2").Value = 24.841114
$ws.Range("I2").Value = 0.2946400644635011
$ws.Range("J2").Value = 0.3116548779253407
$ws.Range("M2").Value = 5.685162000000001
$ws.Range("N2").Value = 17.055486
$ws.Range("O2").Value = 0.06448203205218037
$ws.Range("P2").Value = 0.07835526280970742
$ws.Range("Q2").Value = 47.07525245015601
$ws.Range("R2").Value = 423.6772720514041
$ws.Range("S2").Value = 0.01899899008059197
$ws.Range("T2").Value = 0.02441979986576736
$ws.Range("G3").Value = 8.280371333333333
$ws.Range("H3").Value = 24.841114
$ws.Range("I3").Value = 0.2946400644635011
$ws.Range("J3").Value = 0.3116548779253407
$ws.Range("M3").Value = 5.927452333333334
$ws.Range("O3").Value = 0.06723012841951932
$ws.Range("P3").Value = 0.08169460876758601
$ws.Range("Q3").Value = 49.08150638063312
$ws.Range("R3").Value = 441.7335574256981
$ws.Range("S3").Value = 0.01980868937141663
$ws.Range("T3").Value = 0.02546052332262049
$ws.Range("G4").Value = 8.280371333333333
$ws.Range("H4").Value = 24.841114
$ws.Range("I4").Value = 0.2946400644635011
$ws.Range("J4").Value = 0.3116548779253407
$ws.Range("M4").Value = 13.555466
$ws.Range("N4").Value = 40.666398
$ws.Range("O4").Value = 0.1537483000650186
$ws.Range("P4").Value = 0.1868270597984813
$ws.Range("Q4").Value = 112.2442920763747
$ws.Range("R4").Value = 1010.198628687372
$ws.Range("S4").Value = 0.04530040904231079
$ws.Range("T4").Value = 0.05822556451464601
$ws.Range("G5").Value = 8.280371333333333
$ws.Range("H5").Value = 24.841114
$ws.Range("I5").Value = 0.2946400644635011
$ws.Range("J5").Value = 0.3116548779253407
$ws.Range("M5").Value = 46.83114999999999
$ws.Range("N5").Value = 93.66229999999999
$ws.Range("O5").Value = 0.5311665200288868
$ws.Range("P5").Value = 0.4302975671207292
$ws.Range("Q5").Value = 387.7793119670333
$ws.Range("R5").Value = 2326.6758718022
$ws.Range("S5").Value = 0.1565029377021648
$ws.Range("T5").Value = 0.134104335752582
$ws.Range("G6").Value = 8.280371333333333
$ws.Range("H6").Value = 24.841114
$ws.Range("I6").Value = 0.2946400644635011
$ws.Range("J6").Value = 0.3116548779253407
$ws.Range("M6").Value = 16.167377
$ws.Range("N6").Value = 48.502131
$ws.Range("O6").Value = 0.1833730194343949
$ws.Range("P6").Value = 0.2228255015034961
$ws.Range("Q6").Value = 133.8718850459927
$ws.Range("R6").Value = 1204.846965413934
$ws.Range("S6").Value = 0.05402903826701696
$ws.Range("T6").Value = 0.06944465446972491
$ws.Range("I7").Value = 0.405746032520008
$ws.Range("J7").Value = 0.4291769704298953
$ws.Range("M7").Value = 5.685162000000001
$ws.Range("N7").Value = 17.055486
$ws.Range("O7").Value = 0.06448203205218037
$ws.Range("P7").Value = 0.07835526280970742
$ws.Range("Q7").Value = 64.826882746948
$ws.Range("R7").Value = 583.4419447225321
$ws.Range("S7").Value = 0.02616332867400017
$ws.Range("T7").Value = 0.03362827430990847
$ws.Range("I8").Value = 0.405746032520008
$ws.Range("J8").Value = 0.4291769704298953
$ws.Range("M8").Value = 5.927452333333334
$ws.Range("O8").Value = 0.06723012841951932
$ws.Range("P8").Value = 0.08169460876758601
$ws.Range("Q8").Value = 67.58967596721489
$ws.Range("R8").Value = 608.307083704934
$ws.Range("S8").Value = 0.0272783578720306
$ws.Range("T8").Value = 0.03506144469132812
$ws.Range("I9").Value = 0.405746032520008
$ws.Range("J9").Value = 0.4291769704298953
$ws.Range("M9").Value = 13.555466
$ws.Range("N9").Value = 40.666398
$ws.Range("O9").Value = 0.1537483000650186
$ws.Range("P9").Value = 0.1868270597984813
$ws.Range("Q9").Value = 154.5705478510973
$ws.Range("R9").Value = 1391.134930659876
$ws.Range("S9").Value = 0.06238276275807698
$ws.Range("T9").Value = 0.08018187151863708
$ws.Range("I10").Value = 0.405746032520008
$ws.Range("J10").Value = 0.4291769704298953
$ws.Range("M10").Value = 46.83114999999999
$ws.Range("N10").Value = 93.66229999999999
$ws.Range("O10").Value = 0.5311665200288868
$ws.Range("P10").Value = 0.4302975671207292
$ws.Range("Q10").Value = 534.0072050637665
$ws.Range("R10").Value = 3204.043230382599
$ws.Range("S10").Value = 0.2155187081091802
$ws.Range("T10").Value = 0.1846738062402291
$ws.Range("I11").Value = 0.405746032520008
$ws.Range("J11").Value = 0.4291769704298953
$ws.Range("M11").Value = 16.167377
$ws.Range("N11").Value = 48.502131
$ws.Range("O11").Value = 0.1833730194343949
$ws.Range("P11").Value = 0.2228255015034961
$ws.Range("Q11").Value = 184.3537005813913
$ws.Range("R11").Value = 1659.183305232522
$ws.Range("S11").Value = 0.07440287510672006
$ws.Range("T11").Value = 0.09563157366979254
$ws.Range("G12").Value = 1.864050333333333
$ws.Range("H12").Value = 5.592150999999999
$ws.Range("I12").Value = 0.06632841551025578
$ws.Range("J12").Value = 0.07015873512134246
$ws.Range("M12").Value = 5.685162000000001
$ws.Range("N12").Value = 17.055486
$ws.Range("O12").Value = 0.06448203205218037
$ws.Range("P12").Value = 0.07835526280970742
$ws.Range("Q12").Value = 10.597428121154
$ws.Range("R12").Value = 95.37685309038601
$ws.Range("S12").Value = 0.004276991014902651
$ws.Range("T12").Value = 0.005497306128829439
$ws.Range("G13").Value = 1.864050333333333
$ws.Range("H13").Value = 5.592150999999999
$ws.Range("I13").Value = 0.06632841551025578
$ws.Range("J13").Value = 0.07015873512134246
$ws.Range("M13").Value = 5.927452333333334
$ws.Range("O13").Value = 0.06723012841951932
$ws.Range("P13").Value = 0.08169460876758601
$ws.Range("Q13").Value = 11.04906949776744
$ws.Range("R13").Value = 99.44162547990699
$ws.Range("S13").Value = 0.004459267892617733
$ws.Range("T13").Value = 0.005731590417366769
$ws.Range("G14").Value = 1.864050333333333
$ws.Range("H14").Value = 5.592150999999999
$ws.Range("I14").Value = 0.06632841551025578
$ws.Range("J14").Value = 0.07015873512134246
$ws.Range("M14").Value = 13.555466
$ws.Range("N14").Value = 40.666398
$ws.Range("O14").Value = 0.1537483000650186
$ws.Range("P14").Value = 0.1868270597984813
$ws.Range("Q14").Value = 25.26807091578867
$ws.Range("R14").Value = 227.412638242098
$ws.Range("S14").Value = 0.01019788113070804
$ws.Range("T14").Value = 0.01310755020190086
$ws.Range("G15").Value = 1.864050333333333
$ws.Range("H15").Value = 5.592150999999999
$ws.Range("I15").Value = 0.06632841551025578
$ws.Range("J15").Value = 0.07015873512134246
$ws.Range("M15").Value = 46.83114999999999
$ws.Range("N15").Value = 93.66229999999999
$ws.Range("O15").Value = 0.5311665200288868
$ws.Range("P15").Value = 0.4302975671207292
$ws.Range("Q15").Value = 87.29562076788331
$ws.Range("R15").Value = 523.7737246072999
$ws.Range("S15").Value = 0.0352314336456126
$ws.Range("T15").Value = 0.03018913303498132
$ws.Range("G16").Value = 1.864050333333333
$ws.Range("H16").Value = 5.592150999999999
$ws.Range("I16").Value = 0.06632841551025578
$ws.Range("J16").Value = 0.07015873512134246
$ws.Range("M16").Value = 16.167377
$ws.Range("N16").Value = 48.502131
$ws.Range("O16").Value = 0.1833730194343949
$ws.Range("P16").Value = 0.2228255015034961
$ws.Range("Q16").Value = 30.13680448597566
$ws.Range("R16").Value = 271.231240373781
$ws.Range("S16").Value = 0.01216284182641475
$ws.Range("T16").Value = 0.01563315533826408
$ws.Range("G17").Value = 4.6029105
$ws.Range("H17").Value = 9.205821
$ws.Range("I17").Value = 0.1637851482553954
$ws.Range("J17").Value = 0.1154955860658076
$ws.Range("M17").Value = 5.685162000000001
$ws.Range("N17").Value = 17.055486
$ws.Range("O17").Value = 0.06448203205218037
$ws.Range("P17").Value = 0.07835526280970742
$ws.Range("Q17").Value = 26.16829186400101
$ws.Range("R17").Value = 157.009751184006
$ws.Range("S17").Value = 0.01056119917947552
$ws.Range("T17").Value = 0.009049686999547538
$ws.Range("G18").Value = 4.6029105
$ws.Range("H18").Value = 9.205821
$ws.Range("I18").Value = 0.1637851482553954
$ws.Range("J18").Value = 0.1154955860658076
$ws.Range("M18").Value = 5.927452333333334
$ws.Range("O18").Value = 0.06723012841951932
$ws.Range("P18").Value = 0.08169460876758601
$ws.Range("Q18").Value = 27.2835325833495
$ws.Range("R18").Value = 163.701195500097
$ws.Range("S18").Value = 0.01101129655042024
$ws.Range("T18").Value = 0.009435366718029211
$ws.Range("G19").Value = 4.6029105
$ws.Range("H19").Value = 9.205821
$ws.Range("I19").Value = 0.1637851482553954
$ws.Range("J19").Value = 0.1154955860658076
$ws.Range("M19").Value = 13.555466
$ws.Range("N19").Value = 40.666398
$ws.Range("O19").Value = 0.1537483000650186
$ws.Range("P19").Value = 0.1868270597984813
$ws.Range("Q19").Value = 62.39459678379301
$ws.Range("R19").Value = 374.367580702758
$ws.Range("S19").Value = 0.02518168812016409
$ws.Range("T19").Value = 0.02157770076437728
$ws.Range("G20").Value = 4.6029105
$ws.Range("H20").Value = 9.205821
$ws.Range("I20").Value = 0.1637851482553954
$ws.Range("J20").Value = 0.1154955860658076
$ws.Range("M20").Value = 46.83114999999999
$ws.Range("N20").Value = 93.66229999999999
$ws.Range("O20").Value = 0.5311665200288868
$ws.Range("P20").Value = 0.4302975671207292
$ws.Range("Q20").Value = 215.559592062075
$ws.Range("R20").Value = 862.2383682482999
$ws.Range("S20").Value = 0.08699718723123367
$ws.Range("T20").Value = 0.0496974696972998
$ws.Range("G21").Value = 4.6029105
$ws.Range("H21").Value = 9.205821
$ws.Range("I21").Value = 0.1637851482553954
$ws.Range("J21").Value = 0.1154955860658076
$ws.Range("M21").Value = 16.167377
$ws.Range("N21").Value = 48.502131
$ws.Range("O21").Value = 0.1833730194343949
$ws.Range("P21").Value = 0.2228255015034961
$ws.Range("Q21").Value = 74.4169893507585
$ws.Range("R21").Value = 446.501936104551
$ws.Range("S21").Value = 0.03003377717410188
$ws.Range("T21").Value = 0.02573536188655378
$ws.Range("G22").Value = 1.953192
$ws.Range("H22").Value = 5.859576
$ws.Range("I22").Value = 0.0695003392508397
$ws.Range("J22").Value = 0.0735138304576138
$ws.Range("M22").Value = 5.685162000000001
$ws.Range("N22").Value = 17.055486
$ws.Range("O22").Value = 0.06448203205218037
$ws.Range("P22").Value = 0.07835526280970742
$ws.Range("Q22").Value = 11.104212937104
$ws.Range("R22").Value = 99.93791643393601
$ws.Range("S22").Value = 0.004481523103210055
$ws.Range("T22").Value = 0.005760195505654603
$ws.Range("G23").Value = 1.953192
$ws.Range("H23").Value = 5.859576
$ws.Range("I23").Value = 0.0695003392508397
$ws.Range("J23").Value = 0.0735138304576138
$ws.Range("M23").Value = 5.927452333333334
$ws.Range("O23").Value = 0.06723012841951932
$ws.Range("P23").Value = 0.08169460876758601
$ws.Range("Q23").Value = 11.577452477848
$ws.Range("R23").Value = 104.197072300632
$ws.Range("S23").Value = 0.004672516733034112
$ws.Range("T23").Value = 0.006005683618241407
$ws.Range("G24").Value = 1.953192
$ws.Range("H24").Value = 5.859576
$ws.Range("I24").Value = 0.0695003392508397
$ws.Range("J24").Value = 0.0735138304576138
$ws.Range("M24").Value = 13.555466
$ws.Range("N24").Value = 40.666398
$ws.Range("O24").Value = 0.1537483000650186
$ws.Range("P24").Value = 0.1868270597984813
$ws.Range("Q24").Value = 26.476427747472
$ws.Range("R24").Value = 238.287849727248
$ws.Range("S24").Value = 0.01068555901375869
$ws.Range("T24").Value = 0.01373437279892003
$ws.Range("G25").Value = 1.953192
$ws.Range("H25").Value = 5.859576
$ws.Range("I25").Value = 0.0695003392508397
$ws.Range("J25").Value = 0.0735138304576138
$ws.Range("M25").Value = 46.83114999999999
$ws.Range("N25").Value = 93.66229999999999
$ws.Range("O25").Value = 0.5311665200288868
$ws.Range("P25").Value = 0.4302975671207292
$ws.Range("Q25").Value = 91.47022753079997
$ws.Range("R25").Value = 548.8213651847999
$ws.Range("S25").Value = 0.03691625334069557
$ws.Range("T25").Value = 0.03163282239563698
$ws.Range("G26").Value = 1.953192
$ws.Range("H26").Value = 5.859576
$ws.Range("I26").Value = 0.0695003392508397
$ws.Range("J26").Value = 0.0735138304576138
$ws.Range("M26").Value = 16.167377
$ws.Range("N26").Value = 48.502131
$ws.Range("O26").Value = 0.1833730194343949
$ws.Range("P26").Value = 0.2228255015034961
$ws.Range("Q26").Value = 31.57799141738399
$ws.Range("R26").Value = 284.201922756456
$ws.Range("S26").Value = 0.01216284182641475
$ws.Range("T26").Value = 0.01638075613916078
